$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 4 data rows (rows 2-5, corresponding to years 1984-1987),
# shifting the remaining data up. This also updates the used range /
# dimension from A1:E42 to A1:E38.
$ws.Range("A2:E5").Delete()
